$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary Table")
$wsCooc = $wb.Worksheets.Item("Cooccurrence")
$wsAssoc = $wb.Worksheets.Item("Associations")

# ---------------------------------------------------------------------------
# 1) Summary Table sheet
#    - header text trimmed, header style switched to the existing bold/boxed
#      style (index 2, already present on the Associations sheet header)
#    - the placeholder "---------" row (row 2) is removed, shifting every
#      data row up by one
#    - page margins widened back to the Excel defaults
# ---------------------------------------------------------------------------

$wsSummary.Range("A1").Value = "Keyword"
$wsSummary.Range("B1").Value = "Keyword Count"
$wsSummary.Range("C1").Value = "Short Summary"
$wsSummary.Range("D1").Value = "Source URL"
$wsSummary.Range("E1").Value = "Detailed Summary"

# Reuse the pre-existing bold/bordered style (as seen on Associations!A1)
# by copying formats only, so no new style entries get created.
$wsAssoc.Range("A1").Copy()
$wsSummary.Range("A1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsSummary.Rows(2).Delete()

$ps = $wsSummary.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2) Cooccurrence sheet - populate with the keyword co-occurrence table
# ---------------------------------------------------------------------------

$wsCooc.Cells.Item(1,1).Value = "source"
$wsCooc.Cells.Item(1,2).Value = "target"
$wsCooc.Cells.Item(1,3).Value = "count"

$wsAssoc.Range("A1").Copy()
$wsCooc.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$coocRows = @(
    @("新质生产力", "科技成果转化", 1),
    @("人工智能", "科技成果转化", 1),
    @("知识产权保护", "联合研究中心", 1),
    @("知识产权保护", "量子计算", 1),
    @("知识产权保护", "量子通信", 1),
    @("联合研究中心", "量子计算", 1),
    @("联合研究中心", "量子通信", 1),
    @("量子计算", "量子通信", 1)
)

$r = 2
foreach ($row in $coocRows) {
    $wsCooc.Cells.Item($r, 1).Value = $row[0]
    $wsCooc.Cells.Item($r, 2).Value = $row[1]
    $wsCooc.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Associations sheet - replace placeholder row and append keyword counts
# ---------------------------------------------------------------------------

$assocRows = @(
    @("新质生产力", 1),
    @("科技成果转化", 4),
    @("人工智能", 5),
    @("量子计算", 2),
    @("量子通信", 2),
    @("联合研究中心", 2),
    @("知识产权保护", 2),
    @("工业互联网安全", 1)
)

$r = 2
foreach ($row in $assocRows) {
    $wsAssoc.Cells.Item($r, 1).Value = $row[0]
    $wsAssoc.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
